$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.296.76"
$ws.Range("E2").Value = "  +4.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.488.77"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.16"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.14"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.25"
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.26"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.14"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.479.79"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.234.09"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.54"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.63"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.48"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.14"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  +6.30%  "
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("E31").Value = "  +5.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.60"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.73"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0779"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.86"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.46"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.952.38"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.19"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.80"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +13.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.78"
$ws.Range("E51").Value = "  +3.64%  "
